# Auto-generated: update the "Authors" column (column E) of Sheet1
# for the 34 rows whose author-list text was cleaned up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = '[Ying%Liang%NULL%1,                 Jing-Jin%Liang%NULL%1,                 Qing-Tao%Zhou%NULL%1,                 Xiao-Guang%Li%NULL%1,                 Fei%Lin%NULL%1,                 Zhong-Hua%Deng%NULL%1,                 Bi-Ying%Zhang%NULL%1,                 Lu%Li%NULL%1,                 Xiao-Hua%Wang%NULL%1,                 Hong%Zhu%NULL%1,                 Qing-Bian%Ma%NULL%1,                 Xiao-Mei%Tong%NULL%1,                 Jie%Xu%NULL%1,                 Yong-Chang%Sun%NULL%1,                 Pei-Fang%Wei%NULL%0,                 Pei-Fang%Wei%NULL%0]'
$ws.Cells.Item(3, 5).Value2 = '[Xiao-Wei%Xu%NULL%0,                 Xiao-Xin%Wu%NULL%0,                 Xian-Gao%Jiang%NULL%0,                 Kai-Jin%Xu%NULL%0,                 Ling-Jun%Ying%NULL%0,                 Chun-Lian%Ma%NULL%0,                 Shi-Bo%Li%NULL%0,                 Hua-Ying%Wang%NULL%0,                 Sheng%Zhang%NULL%0,                 Hai-Nv%Gao%NULL%0,                 Ji-Fang%Sheng%NULL%0,                 Hong-Liu%Cai%NULL%0,                 Yun-Qing%Qiu%NULL%0,                 Lan-Juan%Li%NULL%0]'
$ws.Cells.Item(4, 5).Value2 = '[Pavan K.%Bhatraju%NULL%0,                 Bijan J.%Ghassemieh%NULL%0,                 Michelle%Nichols%NULL%0,                 Richard%Kim%NULL%0,                 Keith R.%Jerome%NULL%0,                 Arun K.%Nalla%NULL%0,                 Alexander L.%Greninger%NULL%0,                 Sudhakar%Pipavath%NULL%0,                 Mark M.%Wurfel%NULL%0,                 Laura%Evans%NULL%0,                 Patricia A.%Kritek%NULL%0,                 T. Eoin%West%NULL%0,                 Andrew%Luks%NULL%0,                 Anthony%Gerbino%NULL%0,                 Chris R.%Dale%NULL%0,                 Jason D.%Goldman%NULL%0,                 Shane%O''Mahony%NULL%0,                 Carmen%Mikacenic%NULL%0]'
$ws.Cells.Item(5, 5).Value2 = '[Matt%Arentz%xref no email%0,           Eric%Yim%xref no email%0,           Lindy%Klaff%xref no email%0,           Sharukh%Lokhandwala%xref no email%0,           Francis X.%Riedo%xref no email%0,           Maria%Chong%xref no email%0,           Melissa%Lee%xref no email%0]'
$ws.Cells.Item(6, 5).Value2 = '[NULL%NULL%NULL%0,                 NULL%NULL%NULL%0,                 NULL%NULL%NULL%0,                 Stephanie%Bialek%NULL%1,                 Ellen%Boundy%NULL%1,                 Virginia%Bowen%NULL%1,                 Nancy%Chow%NULL%1,                 Amanda%Cohn%NULL%0,                 Nicole%Dowling%NULL%1,                 Sascha%Ellington%NULL%1,                 Ryan%Gierke%NULL%1,                 Aron%Hall%NULL%1,                 Jessica%MacNeil%NULL%1,                 Priti%Patel%NULL%1,                 Georgina%Peacock%NULL%1,                 Tamara%Pilishvili%NULL%1,                 Hilda%Razzaghi%NULL%1,                 Nia%Reed%NULL%1,                 Matthew%Ritchey%NULL%1,                 Erin%Sauber-Schatz%NULL%1]'
$ws.Cells.Item(7, 5).Value2 = '[Jianlei%Cao%NULL%0,                 Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%0,                 Wenlin%Cheng%NULL%0,                 Lei%Yu%NULL%0,                 Ya-Kun%Liu%NULL%0,                 Xiaoyong%Hu%NULL%0,                 Qiang%Liu%NULL%0]'
$ws.Cells.Item(8, 5).Value2 = '[Jun%Chen%NULL%0,                 Tangkai%Qi%NULL%0,                 Li%Liu%NULL%0,                 Yun%Ling%NULL%0,                 Zhiping%Qian%NULL%0,                 Tao%Li%NULL%0,                 Feng%Li%NULL%0,                 Qingnian%Xu%NULL%0,                 Yuyi%Zhang%NULL%0,                 Shuibao%Xu%NULL%0,                 Zhigang%Song%NULL%0,                 Yigang%Zeng%NULL%0,                 Yinzhong%Shen%NULL%0,                 Yuxin%Shi%NULL%0,                 Tongyu%Zhu%NULL%0,                 Hongzhou%Lu%NULL%0]'
$ws.Cells.Item(9, 5).Value2 = '[Na%Zhu%NULL%0,                 Dingyu%Zhang%NULL%0,                 Wenling%Wang%NULL%0,                 Xingwang%Li%NULL%0,                 Bo%Yang%NULL%0,                 Jingdong%Song%NULL%0,                 Xiang%Zhao%NULL%0,                 Baoying%Huang%NULL%0,                 Weifeng%Shi%NULL%0,                 Roujian%Lu%NULL%0,                 Peihua%Niu%NULL%0,                 Faxian%Zhan%NULL%0,                 Xuejun%Ma%NULL%0,                 Dayan%Wang%NULL%0,                 Wenbo%Xu%NULL%0,                 Guizhen%Wu%NULL%0,                 George F.%Gao%NULL%0,                 Wenjie%Tan%NULL%0]'
$ws.Cells.Item(10, 5).Value2 = '[Christopher M%Petrilli%NULL%5,                 Simon A%Jones%NULL%1,                 Jie%Yang%NULL%1,                 Harish%Rajagopalan%NULL%1,                 Luke%O''Donnell%NULL%1,                 Yelena%Chernyak%NULL%1,                 Katie A%Tobin%NULL%1,                 Robert J%Cerfolio%NULL%1,                 Fritz%Francois%NULL%2,                 Leora I%Horwitz%NULL%1]'
$ws.Cells.Item(11, 5).Value2 = '[Barnaby%D.P.%coreGivesNoEmail%0,               Becker%L.B.%coreGivesNoEmail%0,               Chelico%J.D.%coreGivesNoEmail%0,               Cohen%S.L.%coreGivesNoEmail%0,               Cookingham%J.%coreGivesNoEmail%0,               Coppa%K.%coreGivesNoEmail%0,               Crawford%J.M.%coreGivesNoEmail%0,               Davidson%K.W.%coreGivesNoEmail%0,               Diefenbach%M.A.%coreGivesNoEmail%0,               Dominello%A.J.%coreGivesNoEmail%0,               Duer-Hefele%J.%coreGivesNoEmail%0,               Falzon%L.%coreGivesNoEmail%0,               Gitlin%J.%coreGivesNoEmail%0,               Hajizadeh%N.%coreGivesNoEmail%0,               Harvin%T.G.%coreGivesNoEmail%0,               Hirsch%J.S.%coreGivesNoEmail%0,               Hirschwerk%D.A.%coreGivesNoEmail%0,               Kim%E.J.%coreGivesNoEmail%0,               Kozel%Z.M.%coreGivesNoEmail%0,               Marrast%L.M.%coreGivesNoEmail%0,               McGinn%T.%coreGivesNoEmail%0,               Mogavero%J.N.%coreGivesNoEmail%0,               Narasimhan%M.%coreGivesNoEmail%0,               Osorio%G.A.%coreGivesNoEmail%0,               Qiu%M.%coreGivesNoEmail%0,               Richardson%S.%coreGivesNoEmail%0,               Zanos%T.P.%coreGivesNoEmail%0]'
$ws.Cells.Item(12, 5).Value2 = '[Arthur%Simonnet%NULL%0,                 Mikael%Chetboun%NULL%0,                 Julien%Poissy%NULL%0,                 Violeta%Raverdy%NULL%0,                 Jerome%Noulette%NULL%0,                 Jerome%Noulette%NULL%0,                 Alain%Duhamel%NULL%0,                 Julien%Labreuche%NULL%0,                 Daniel%Mathieu%NULL%0,                 Francois%Pattou%francois.pattou@univ-lille.fr%0,                 Merce%Jourdain%mercedes.jourdain@univ-lille.fr%0,                 Merce%Jourdain%mercedes.jourdain@univ-lille.fr%0,                 Robert%Caizzo%NULL%0,                 Robert%Caizzo%NULL%0,                 Morgan%Caplan%NULL%0,                 Nicolas%Cousin%NULL%0,                 Thibault%Duburcq%NULL%0,                 Arthur%Durand%NULL%0,                 Ahmed%El kalioubie%NULL%0,                 Raphael%Favory%NULL%0,                 Bruno%Garcia%NULL%0,                 Patrick%Girardie%NULL%0,                 Julien%Goutay%NULL%0,                 Marion%Houard%NULL%0,                 Emmanuelle%Jaillette%NULL%0,                 Nicolas%Kostuj%NULL%0,                 Geoffrey%Ledoux%NULL%0,                 Daniel%Mathieu%NULL%0,                 Anne Sophie%Moreau%NULL%0,                 Christopher%Niles%NULL%0,                 Saad%Nseir%NULL%0,                 Thierry%Onimus%NULL%0,                 Erika%Parmentier%NULL%0,                 Sebastien%Préau%NULL%0,                 Laurent%Robriquet%NULL%0,                 Anahita%Rouze%NULL%0,                 Sophie%Six%NULL%0,                 Hélène%Verkindt%NULL%0]'
$ws.Cells.Item(13, 5).Value2 = '[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,               Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%0,               Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%0,               Subert-Salas%Lizandra%coreGivesNoEmail%0]'
$ws.Cells.Item(14, 5).Value2 = '[Chaomin%Wu%xref no email%0,           Xiaoyan%Chen%xref no email%0,           Yanping%Cai%xref no email%0,           Jia’an%Xia%xref no email%0,           Xing%Zhou%xref no email%0,           Sha%Xu%xref no email%0,           Hanping%Huang%xref no email%0,           Li%Zhang%xref no email%0,           Xia%Zhou%xref no email%0,           Chunling%Du%xref no email%0,           Yuye%Zhang%xref no email%0,           Juan%Song%xref no email%0,           Sijiao%Wang%xref no email%0,           Yencheng%Chao%xref no email%0,           Zeyong%Yang%xref no email%0,           Jie%Xu%xref no email%0,           Xin%Zhou%xref no email%0,           Dechang%Chen%xref no email%0,           Weining%Xiong%xref no email%0,           Lei%Xu%xref no email%0,           Feng%Zhou%xref no email%0,           Jinjun%Jiang%xref no email%0,           Chunxue%Bai%xref no email%0,           Junhua%Zheng%xref no email%0,           Yuanlin%Song%xref no email%0]'
$ws.Cells.Item(15, 5).Value2 = '[Cheng-Yu%Chen%NULL%1,                 Chen-Hsen%Lee%NULL%1,                 Cheng-Yi%Liu%NULL%1,                 Jia-Horng%Wang%NULL%1,                 Lee-Min%Wang%NULL%1,                 Reury-Perng%Perng%NULL%1]'
$ws.Cells.Item(16, 5).Value2 = '[Jiqian%Xu%NULL%2,                 Xiaobo%Yang%NULL%1,                 Luyu%Yang%NULL%1,                 Xiaojing%Zou%NULL%0,                 Yaxin%Wang%NULL%0,                 Yongran%Wu%NULL%0,                 Ting%Zhou%NULL%1,                 Yin%Yuan%NULL%1,                 Hong%Qi%NULL%1,                 Shouzhi%Fu%NULL%0,                 Hong%Liu%NULL%0,                 Jia''an%Xia%NULL%1,                 Zhengqin%Xu%NULL%1,                 Yuan%Yu%NULL%0,                 Ruiting%Li%NULL%0,                 Yaqi%Ouyang%NULL%1,                 Rui%Wang%NULL%0,                 Lehao%Ren%NULL%1,                 Yingying%Hu%NULL%1,                 Dan%Xu%NULL%0,                 Xin%Zhao%NULL%1,                 Shiying%Yuan%shiying_yuan@163.com%0,                 Dingyu%Zhang%1813886398@qq.com%0,                 You%Shang%you_shanghust@163.com%0]'
$ws.Cells.Item(17, 5).Value2 = '[Barnaby Edward%Young%xref no email%0,           Sean Wei Xiang%Ong%xref no email%0,           Shirin%Kalimuddin%xref no email%0,           Jenny G.%Low%xref no email%0,           Seow Yen%Tan%xref no email%0,           Jiashen%Loh%xref no email%0,           Oon-Tek%Ng%xref no email%0,           Kalisvar%Marimuthu%xref no email%0,           Li Wei%Ang%xref no email%0,           Tze Minn%Mak%xref no email%0,           Sok Kiang%Lau%xref no email%0,           Danielle E.%Anderson%xref no email%0,           Kian Sing%Chan%xref no email%0,           Thean Yen%Tan%xref no email%0,           Tong Yong%Ng%xref no email%0,           Lin%Cui%xref no email%0,           Zubaidah%Said%xref no email%0,           Lalitha%Kurupatham%xref no email%0,           Mark I-Cheng%Chen%xref no email%0,           Monica%Chan%xref no email%0,           Shawn%Vasoo%xref no email%0,           Lin-Fa%Wang%xref no email%0,           Boon Huan%Tan%xref no email%0,           Raymond Tzer Pin%Lin%xref no email%0,           Vernon Jian Ming%Lee%xref no email%0,           Yee-Sin%Leo%xref no email%0,           David Chien%Lye%xref no email%0]'
$ws.Cells.Item(18, 5).Value2 = '[Wei-jie%Guan%NULL%0,                 Wen-hua%Liang%NULL%0,                 Yi%Zhao%NULL%1,                 Heng-rui%Liang%NULL%1,                 Zi-sheng%Chen%NULL%1,                 Yi-min%Li%NULL%1,                 Xiao-qing%Liu%NULL%1,                 Ru-chong%Chen%NULL%0,                 Chun-li%Tang%NULL%0,                 Tao%Wang%NULL%0,                 Chun-quan%Ou%NULL%0,                 Li%Li%NULL%0,                 Ping-yan%Chen%NULL%0,                 Ling%Sang%NULL%1,                 Wei%Wang%NULL%0,                 Jian-fu%Li%NULL%1,                 Cai-chen%Li%NULL%1,                 Li-min%Ou%NULL%1,                 Bo%Cheng%NULL%1,                 Shan%Xiong%NULL%1,                 Zheng-yi%Ni%NULL%0,                 Jie%Xiang%NULL%0,                 Yu%Hu%NULL%0,                 Lei%Liu%NULL%0,                 Hong%Shan%NULL%0,                 Chun-liang%Lei%NULL%0,                 Yi-xiang%Peng%NULL%0,                 Li%Wei%NULL%0,                 Yong%Liu%NULL%0,                 Ya-hua%Hu%NULL%0,                 Peng%Peng%NULL%0,                 Jian-ming%Wang%NULL%0,                 Ji-yang%Liu%NULL%0,                 Zhong%Chen%NULL%0,                 Gang%Li%NULL%0,                 Zhi-jian%Zheng%NULL%0,                 Shao-qin%Qiu%NULL%0,                 Jie%Luo%NULL%0,                 Chang-jiang%Ye%NULL%0,                 Shao-yong%Zhu%NULL%0,                 Lin-ling%Cheng%NULL%1,                 Feng%Ye%NULL%1,                 Shi-yue%Li%NULL%0,                 Jin-ping%Zheng%NULL%1,                 Nuo-fu%Zhang%NULL%1,                 Nan-shan%Zhong%NULL%0,                 Jian-xing%He%NULL%0]'
$ws.Cells.Item(19, 5).Value2 = '[Fei%Zhou%NULL%0,                 Ting%Yu%NULL%0,                 Ronghui%Du%NULL%0,                 Guohui%Fan%NULL%0,                 Ying%Liu%NULL%0,                 Zhibo%Liu%NULL%0,                 Jie%Xiang%NULL%0,                 Yeming%Wang%NULL%0,                 Bin%Song%NULL%0,                 Xiaoying%Gu%NULL%0,                 Lulu%Guan%NULL%0,                 Yuan%Wei%NULL%0,                 Hui%Li%NULL%0,                 Xudong%Wu%NULL%0,                 Jiuyang%Xu%NULL%0,                 Shengjin%Tu%NULL%0,                 Yi%Zhang%NULL%0,                 Hua%Chen%NULL%0,                 Bin%Cao%NULL%0]'
$ws.Cells.Item(20, 5).Value2 = '[Corrado%Lodigiani%NULL%1,                 Giacomo%Iapichino%NULL%1,                 Luca%Carenzo%NULL%1,                 Maurizio%Cecconi%NULL%1,                 Paola%Ferrazzi%NULL%1,                 Tim%Sebastian%NULL%1,                 Nils%Kucher%NULL%1,                 Jan-Dirk%Studt%NULL%1,                 Clara%Sacco%NULL%1,                 Alexia%Bertuzzi%NULL%1,                 Maria Teresa%Sandri%NULL%1,                 Stefano%Barco%NULL%1,                 NULL%NULL%NULL%0]'
$ws.Cells.Item(21, 5).Value2 = '[F.A.%Klok%NULL%1,                 M.J.H.A.%Kruip%NULL%1,                 N.J.M.%van der Meer%NULL%1,                 M.S.%Arbous%NULL%1,                 D.A.M.P.J.%Gommers%NULL%1,                 K.M.%Kant%NULL%1,                 F.H.J.%Kaptein%NULL%1,                 J.%van Paassen%NULL%1,                 M.A.M.%Stals%NULL%1,                 M.V.%Huisman%NULL%1,                 H.%Endeman%NULL%1]'
$ws.Cells.Item(22, 5).Value2 = '[Shaoqing%Lei%NULL%3,                 Fang%Jiang%NULL%1,                 Wating%Su%NULL%1,                 Chang%Chen%NULL%1,                 Jingli%Chen%NULL%1,                 Wei%Mei%NULL%1,                 Li-Ying%Zhan%NULL%1,                 Yifan%Jia%NULL%1,                 Liangqing%Zhang%NULL%1,                 Danyong%Liu%NULL%1,                 Zhong-Yuan%Xia%xiazhongyuan2005@aliyun.com%1,                 Zhengyuan%Xia%zyxia@hkucc.hku.hk%1]'
$ws.Cells.Item(23, 5).Value2 = '[Annemarie B%Docherty%NULL%1,                 Ewen M%Harrison%NULL%1,                 Christopher A%Green%NULL%1,                 Hayley E%Hardwick%NULL%1,                 Riinu%Pius%NULL%1,                 Lisa%Norman%NULL%1,                 Karl A%Holden%NULL%1,                 Jonathan M%Read%NULL%1,                 Frank%Dondelinger%NULL%1,                 Gail%Carson%NULL%1,                 Laura%Merson%NULL%1,                 James%Lee%NULL%1,                 Daniel%Plotkin%NULL%1,                 Louise%Sigfrid%NULL%1,                 Sophie%Halpin%NULL%1,                 Clare%Jackson%NULL%1,                 Carrol%Gamble%NULL%1,                 Peter W%Horby%NULL%1,                 Jonathan S%Nguyen-Van-Tam%NULL%1,                 Antonia%Ho%NULL%1,                 Clark D%Russell%NULL%1,                 Jake%Dunning%NULL%1,                 Peter JM%Openshaw%NULL%1,                 J Kenneth%Baillie%NULL%1,                 Malcolm G%Semple%NULL%1]'
$ws.Cells.Item(24, 5).Value2 = '[Rong-Hui%Du%NULL%0,                 Li-Min%Liu%NULL%0,                 Wen%Yin%NULL%0,                 Wen%Wang%NULL%0,                 Lu-Lu%Guan%NULL%0,                 Ming-Li%Yuan%NULL%0,                 Yu-Lei%Li%NULL%0,                 Yi%Hu%NULL%0,                 Xu-Yan%Li%NULL%0,                 Bing%Sun%NULL%0,                 Peng%Peng%NULL%0,                 Huan-Zhong%Shi%NULL%0]'
$ws.Cells.Item(25, 5).Value2 = '[
Ming%Hu%NULL%1,                 Qiang%Zhou%NULL%1,                 Ruiqiang%Zheng%NULL%1,                 Xuyan%Li%NULL%1,                 Jianmin%Ling%NULL%1,                 Yumei%Chen%NULL%1,                 Jing%Jia%NULL%1,                 Cuihong%Xie%xiecuihong08@163.com%1]'
$ws.Cells.Item(26, 5).Value2 = '[Evgeny V.%Fominskiy%NULL%1,                 Anna Mara%Scandroglio%NULL%1,                 Giacomo%Monti%NULL%1,                 Maria Grazia%Calabrò%NULL%1,                 Giovanni%Landoni%NULL%1,                 Antonio%Dell''Acqua%NULL%1,                 Luigi%Beretta%NULL%1,                 Elena%Moizo%NULL%1,                 Alfredo%Ravizza%NULL%1,                 Fabrizio%Monaco%NULL%1,                 Corrado%Campochiaro%NULL%1,                 Marina%Pieri%NULL%1,                 Maria Luisa%Azzolini%NULL%1,                 Giovanni%Borghi%NULL%1,                 Martina%Crivellari%NULL%1,                 Caterina%Conte%NULL%1,                 Cristina%Mattioli%NULL%1,                 Paolo%Silvani%NULL%1,                 Milena%Mucci%NULL%1,                 Stefano%Turi%NULL%1,                 Stefano%Tentori%NULL%1,                 Martina%Baiardo Redaelli%NULL%1,                 Marianna%Sartorelli%NULL%1,                 Piera%Angelillo%NULL%1,                 Alessandro%Belletti%NULL%1,                 Pasquale%Nardelli%NULL%1,                 Francesco Giuseppe%Nisi%NULL%1,                 Gabriele%Valsecchi%NULL%1,                 Cristina%Barberio%NULL%1,                 Fabio%Ciceri%NULL%1,                 Ary%Serpa Neto%NULL%1,                 Lorenzo%Dagna%NULL%1,                 Rinaldo%Bellomo%NULL%1,                 Alberto%Zangrillo%NULL%1]'
$ws.Cells.Item(27, 5).Value2 = '[Antonelli%Massimo%coreGivesNoEmail%0,               Cabrini%Luca%coreGivesNoEmail%0,               Castelli%Antonio%coreGivesNoEmail%0,               Cecconi%Maurizio%coreGivesNoEmail%0,               Cereda%Danilo%coreGivesNoEmail%0,               Coluccello%Antonio%coreGivesNoEmail%0,               Foti%Giuseppe%coreGivesNoEmail%0,               Fumagalli%Roberto%coreGivesNoEmail%0,               Grasselli%Giacomo%coreGivesNoEmail%0,               Iotti%Giorgio%coreGivesNoEmail%0,               Latronico%Nicola%coreGivesNoEmail%0,               Lorini%Luca%coreGivesNoEmail%0,               Merler%Stefano%coreGivesNoEmail%0,               Natalini%Giuseppe%coreGivesNoEmail%0,               Pesenti%Antonio%coreGivesNoEmail%0,               Piatti%Alessandra%coreGivesNoEmail%0,               Ranieri%Marco Vito%coreGivesNoEmail%0,               Scandroglio%Anna Mara%coreGivesNoEmail%0,               Storti%Enrico%coreGivesNoEmail%0,               Zanella%Alberto%coreGivesNoEmail%0,               Zangrillo%Alberto%coreGivesNoEmail%0]'
$ws.Cells.Item(28, 5).Value2 = '[Nanshan%Chen%NULL%0,                 Min%Zhou%NULL%0,                 Xuan%Dong%NULL%0,                 Jieming%Qu%NULL%0,                 Fengyun%Gong%NULL%0,                 Yang%Han%NULL%0,                 Yang%Qiu%NULL%0,                 Jingli%Wang%NULL%0,                 Ying%Liu%NULL%0,                 Yuan%Wei%NULL%0,                 Jia''an%Xia%NULL%0,                 Ting%Yu%NULL%0,                 Xinxin%Zhang%NULL%0,                 Li%Zhang%NULL%0]'
$ws.Cells.Item(30, 5).Value2 = '[Kin Wing%Choi%xref no email%1,           Tai Nin%Chau%xref no email%1,           Owen%Tsang%xref no email%1,           Eugene%Tso%xref no email%1,           Ming Chee%Chiu%xref no email%1,           Wing Lok%Tong%xref no email%1,           Po Oi%Lee%xref no email%1,           Tak Keung%Ng%xref no email%1,           Wai Fu%Ng%xref no email%1,           Kam Cheong%Lee%xref no email%1,           William%Lam%xref no email%1,           Wai Cho%Yu%xref no email%1,           Jak Yiu%Lai%xref no email%1,           Sik To%Lai%xref no email%1]'
$ws.Cells.Item(31, 5).Value2 = '[Fahad%Al-Hameed%xref no email%1,           Ali S.%Wahla%xref no email%1,           Shakeel%Siddiqui%xref no email%1,           Alaa%Ghabashi%xref no email%2,           Majid%Al-Shomrani%xref no email%1,           Abdulhakeem%Al-Thaqafi%xref no email%1,           Yasir%Tashkandi%xref no email%1]'
$ws.Cells.Item(32, 5).Value2 = '[Yaseen M.%Arabi%xref no email%1,           Ahmed A.%Arifi%xref no email%1,           Hanan H.%Balkhy%xref no email%1,           Hani%Najm%xref no email%1,           Abdulaziz S.%Aldawood%xref no email%1,           Alaa%Ghabashi%xref no email%0,           Hassan%Hawa%xref no email%1,           Adel%Alothman%xref no email%1,           Abdulaziz%Khaldi%xref no email%1,           Basel%Al Raiy%xref no email%1]'
$ws.Cells.Item(33, 5).Value2 = '[Mustafa%Saad%NULL%1,                 Ali S.%Omrani%NULL%1,                 Kamran%Baig%NULL%1,                 Abdelkarim%Bahloul%NULL%1,                 Fatehi%Elzein%NULL%1,                 Mohammad Abdul%Matin%NULL%1,                 Mohei A.A.%Selim%NULL%1,                 Mohammed Al%Mutairi%NULL%1,                 Daifullah Al%Nakhli%NULL%1,                 Amal Y. Al%Aidaroos%NULL%1,                 Nisreen Al%Sherbeeni%NULL%1,                 Hesham I.%Al-Khashan%NULL%1,                 Ziad A.%Memish%zmemish@yahoo.com%1,                 Ali M.%Albarrak%NULL%1]'
$ws.Cells.Item(34, 5).Value2 = '[Ashraf Abdel%Halim%NULL%1,                 Badr%Alsayed%NULL%1,                 Sameh%Embarak%NULL%1,                 Taha%Yaseen%NULL%1,                 Salwa%Dabbous%NULL%1]'
$ws.Cells.Item(35, 5).Value2 = '[Musa A.%Garbati%NULL%1,                 Shamsudeen F.%Fagbo%NULL%1,                 Vicky J.%Fang%NULL%1,                 Leila%Skakni%NULL%1,                 Mercy%Joseph%NULL%1,                 Tariq A.%Wani%NULL%1,                 Benjamin J.%Cowling%NULL%0,                 Malik%Peiris%NULL%1,                 Ahmed%Hakawi%NULL%1,                 Renee W.Y.%Chan%NULL%2,                 Renee W.Y.%Chan%NULL%0]'
$ws.Cells.Item(36, 5).Value2 = '[Mohammed%Al Ghamdi%dr_mohd1391@yahoo.com%1,                 Khalid M.%Alghamdi%K.m.a02@hotmail.com%1,                 Yasmeen%Ghandoora%Yass.m.gh@hotmail.com%1,                 Ameera%Alzahrani%Ameera.alzahrani@outlook.com%1,                 Fatmah%Salah%Fs4l_med2011@yahoo.com%1,                 Abdulmoatani%Alsulami%moatani@hotmail.com%1,                 Mayada F.%Bawayan%mayada-f-b@hotmail.com%1,                 Dhananjay%Vaidya%dvaidya@jhmi.edu%1,                 Trish M.%Perl%tperl@jhmi.edu%1,                 Geeta%Sood%gsood1@jhmi.edu%1]'

Write-Host "Updated Authors (column E) for" 34 "rows"
